# Economic Dashboard weekly refresh (2026-01-23)
# - Rolls each indicator 5-period window forward by one release
# - Re-highlights the date cell of whichever release is the newest this week
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Toggle the "updated this week" yellow highlight on release-date cells ---
# Newly-updated release dates: paint with the existing yellow-fill style (copied from N13)
$highlightOnTargets = @(
  "C11", "C12", "C13", "C14", "C15", "C16", "C19", "C20", "C21", "C22", "C23", "N24", "N25", "N26", "N27", "N51"
)
foreach ($ref in $highlightOnTargets) {
  $ws.Range("N13").Copy() | Out-Null
  $ws.Range($ref).PasteSpecial($xlPasteFormats) | Out-Null
}

# Release dates no longer the newest: revert to the plain (unfilled) style (copied from C3)
$highlightOffTargets = @(
  "C32", "C33", "C34"
)
foreach ($ref in $highlightOffTargets) {
  $ws.Range("C3").Copy() | Out-Null
  $ws.Range($ref).PasteSpecial($xlPasteFormats) | Out-Null
}

$excel.CutCopyMode = $false

# --- Update cell values: new release dates + refreshed 5-period data windows ---
$ws.Range("C11").Value = 45962
$ws.Range("F11").Value = 0.006303109534036899
$ws.Range("G11").Value = 0.003232001498899351
$ws.Range("H11").Value = -0.00461581499440511
$ws.Range("I11").Value = 0.001494209936496294
$ws.Range("J11").Value = 0.01310374189885977
$ws.Range("C12").Value = 45962
$ws.Range("F12").Value = 0.01107097621616567
$ws.Range("G12").Value = 0.02185114503816803
$ws.Range("H12").Value = 0.02290259211345895
$ws.Range("I12").Value = 0.03914728682170551
$ws.Range("J12").Value = 0.0341397460041527
$ws.Range("C13").Value = 45962
$ws.Range("F13").Value = 0.005120414134593743
$ws.Range("G13").Value = 0.004379892056853851
$ws.Range("H13").Value = -0.00357585313661446
$ws.Range("I13").Value = 0.005207743688441013
$ws.Range("J13").Value = 0.002952197115930311
$ws.Range("C14").Value = 45962
$ws.Range("F14").Value = 0.03287172222382833
$ws.Range("G14").Value = 0.02743170978465099
$ws.Range("H14").Value = 0.02153393181883783
$ws.Range("I14").Value = 0.03793325150505573
$ws.Range("J14").Value = 0.03053813621117101
$ws.Range("C15").Value = 45962
$ws.Range("F15").Value = 0.001899602810321532
$ws.Range("G15").Value = 0.002642778380250288
$ws.Range("H15").Value = 0.003374053619106698
$ws.Range("I15").Value = 0.002493308400249195
$ws.Range("J15").Value = 0.004234479711318961
$ws.Range("C16").Value = 45962
$ws.Range("F16").Value = 0.02567156402073092
$ws.Range("G16").Value = 0.02575027269930347
$ws.Range("H16").Value = 0.0245940671714955
$ws.Range("I16").Value = 0.02312614601653996
$ws.Range("J16").Value = 0.02374204688350444
$ws.Range("C19").Value = 45962
$ws.Range("F19").Value = 0.0006937237425563847
$ws.Range("G19").Value = -0.001064425459726515
$ws.Range("H19").Value = 0.0005713461581133839
$ws.Range("I19").Value = 0.001005025125628167
$ws.Range("J19").Value = 0.00263888966212189
$ws.Range("C20").Value = 45962
$ws.Range("F20").Value = 0.01007203916779632
$ws.Range("G20").Value = 0.01168973358412177
$ws.Range("H20").Value = 0.01508168307081075
$ws.Range("I20").Value = 0.01547352826862074
$ws.Range("J20").Value = 0.01500856666215786
$ws.Range("C21").Value = 45962
$ws.Range("F21").Value = 3.5
$ws.Range("G21").Value = 3.7
$ws.Range("H21").Value = 4
$ws.Range("I21").Value = 4.1
$ws.Range("J21").Value = 4.3
$ws.Range("C22").Value = 45992
$ws.Range("F22").Value = 16.481
$ws.Range("G22").Value = 16.117
$ws.Range("H22").Value = 15.807
$ws.Range("I22").Value = 16.663
$ws.Range("J22").Value = 16.916
$ws.Range("C23").Value = 45992
$ws.Range("F23").Value = -0.04761629586824606
$ws.Range("G23").Value = -0.05221993531314308
$ws.Range("H23").Value = -0.04524039623097369
$ws.Range("I23").Value = 0.02114229685010411
$ws.Range("J23").Value = 0.05976694649793263
$ws.Range("N24").Value = 45962
$ws.Range("Q24").Value = 0.002073098225740644
$ws.Range("R24").Value = 0.001590597453477116
$ws.Range("S24").Value = 0.002608155986582039
$ws.Range("T24").Value = 0.002622873345935917
$ws.Range("U24").Value = 0.001712126113473822
$ws.Range("N25").Value = 45962
$ws.Range("Q25").Value = 0.02772852363262907
$ws.Range("R25").Value = 0.02678040708789181
$ws.Range("S25").Value = 0.02787442414870654
$ws.Range("T25").Value = 0.02747620854151709
$ws.Range("U25").Value = 0.0260554729423934
$ws.Range("N26").Value = 45962
$ws.Range("Q26").Value = 0.001603546667924283
$ws.Range("R26").Value = 0.002079493359799622
$ws.Range("S26").Value = 0.001894029073346237
$ws.Range("T26").Value = 0.002246302301668779
$ws.Range("U26").Value = 0.002450028147572558
$ws.Range("N27").Value = 45962
$ws.Range("Q27").Value = 0.02791177941627268
$ws.Range("R27").Value = 0.02734349764196662
$ws.Range("S27").Value = 0.02825069249833962
$ws.Range("T27").Value = 0.02912416347215904
$ws.Range("U27").Value = 0.02863047245567936
$ws.Range("N29").Value = 46044
$ws.Range("Q29").Value = 2.2
$ws.Range("S29").Value = 2.26
$ws.Range("T29").Value = 2.27
$ws.Range("N30").Value = 46044
$ws.Range("Q30").Value = 2.31
$ws.Range("R30").Value = 2.34
$ws.Range("T30").Value = 2.33
$ws.Range("C32").Value = 45992
$ws.Range("C33").Value = 45992
$ws.Range("C34").Value = 45992
$ws.Range("N34").Value = 45962
$ws.Range("Q34").Value = 0.008267996311731092
$ws.Range("R34").Value = 0.01042628043129703
$ws.Range("S34").Value = 0.009230301788676142
$ws.Range("T34").Value = 0.01055343698833223
$ws.Range("U34").Value = 0.01240099690797541
$ws.Range("N35").Value = 45962
$ws.Range("Q35").Value = 0.0003711201807323761
$ws.Range("R35").Value = 0.002770616347913979
$ws.Range("S35").Value = -0.0006927352283296884
$ws.Range("T35").Value = 0.001490702319483894
$ws.Range("U35").Value = 0.001590026251228505
$ws.Range("N36").Value = 45962
$ws.Range("Q36").Value = 0.008267996311731092
$ws.Range("R36").Value = 0.01042628043129703
$ws.Range("S36").Value = 0.009230301788676142
$ws.Range("T36").Value = 0.01055343698833223
$ws.Range("U36").Value = 0.01240099690797541
$ws.Range("N47").Value = 46043
$ws.Range("N48").Value = 46043
$ws.Range("R48").Value = 3.6
$ws.Range("S48").Value = 3.59
$ws.Range("T48").Value = 3.56
$ws.Range("U48").Value = 3.51
$ws.Range("N49").Value = 46043
$ws.Range("Q49").Value = 3.83
$ws.Range("R49").Value = 3.86
$ws.Range("S49").Value = 3.82
$ws.Range("T49").Value = 3.77
$ws.Range("U49").Value = 3.72
$ws.Range("N50").Value = 46043
$ws.Range("Q50").Value = 4.26
$ws.Range("R50").Value = 4.3
$ws.Range("S50").Value = 4.24
$ws.Range("T50").Value = 4.17
$ws.Range("U50").Value = 4.15
$ws.Range("N51").Value = 46041
$ws.Range("Q51").Value = 6.09
$ws.Range("R51").Value = 6.06
$ws.Range("S51").Value = 6.16
$ws.Range("T51").Value = 6.15
$ws.Range("U51").Value = 6.18
$ws.Range("N52").Value = 46043
$ws.Range("Q52").Value = 5.88
$ws.Range("R52").Value = 5.95
$ws.Range("S52").Value = 5.87
$ws.Range("T52").Value = 5.82
$ws.Range("U52").Value = 5.83
